$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    if ($val -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        # Numeric-looking text: force Text format so Excel keeps it
        # as a literal string instead of parsing it into a Double.
        $cell.NumberFormat = "@"
        $cell.Value2 = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value2 = $val
    }
}

Set-TextCell 2 4 "24.802.09"
Set-TextCell 2 5 "  +0.74%  "
Set-TextCell 3 4 "1.703.46"
Set-TextCell 3 5 "  +0.68%  "
Set-TextCell 4 5 "  -0.11%  "
Set-TextCell 5 4 "314.42"
Set-TextCell 5 5 "  -0.02%  "
Set-TextCell 6 5 "  -0.26%  "
Set-TextCell 7 4 "0.3993"
Set-TextCell 7 5 "  +2.31%  "
Set-TextCell 8 4 "0.4070"
Set-TextCell 8 5 "  +0.79%  "
Set-TextCell 9 4 "1.002"
Set-TextCell 9 5 "  -0.05%  "
Set-TextCell 10 4 "1.470"
Set-TextCell 10 5 "  -2.13%  "
Set-TextCell 11 4 "53.59"
Set-TextCell 11 5 "  +1.10%  "
Set-TextCell 12 4 "0.08818"
Set-TextCell 12 5 "  +0.75%  "
Set-TextCell 13 4 "26.32"
Set-TextCell 13 5 "  +3.39%  "
Set-TextCell 14 4 "7.509"
Set-TextCell 14 5 "  -0.23%  "
Set-TextCell 15 4 "7.980"
Set-TextCell 15 5 "  +0.21%  "
Set-TextCell 16 5 "  -0.86%  "
Set-TextCell 17 4 "1.738.51"
Set-TextCell 17 5 "  +2.72%  "
Set-TextCell 18 4 "95.44"
Set-TextCell 18 5 "  -3.37%  "
Set-TextCell 19 4 "0.07172"
Set-TextCell 19 5 "  +0.95%  "
Set-TextCell 20 5 "  +4.85%  "
Set-TextCell 21 4 "7.300"
Set-TextCell 21 5 "  -0.13%  "
Set-TextCell 22 4 "1.001"
Set-TextCell 22 5 "  -0.20%  "
Set-TextCell 23 4 "14.45"
Set-TextCell 23 5 "  +1.20%  "
Set-TextCell 24 4 "24.771.74"
Set-TextCell 24 5 "  +0.68%  "
Set-TextCell 25 4 "2.378"
Set-TextCell 25 5 "  +1.00%  "
Set-TextCell 26 4 "2.888"
Set-TextCell 26 5 "  -4.27%  "
Set-TextCell 27 4 "23.12"
Set-TextCell 27 5 "  +1.22%  "
Set-TextCell 28 4 "6.083"
Set-TextCell 28 5 "  +16.58%  "
Set-TextCell 29 4 "161.79"
Set-TextCell 29 5 "  -0.21%  "
Set-TextCell 30 4 "143.95"
Set-TextCell 30 5 "  +4.96%  "
Set-TextCell 31 4 "8.246"
Set-TextCell 31 5 "  -5.06%  "
Set-TextCell 32 4 "2.269"
Set-TextCell 32 5 "  +14.97%  "
Set-TextCell 33 4 "1.920.39"
Set-TextCell 33 5 "  +2.35%  "
Set-TextCell 34 2 "VeChain"
Set-TextCell 34 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell 34 4 "0.03194"
Set-TextCell 34 5 "  +9.02%  "
Set-TextCell 35 2 "Hedera"
Set-TextCell 35 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell 35 4 "0.08585"
Set-TextCell 35 5 "  -2.64%  "
Set-TextCell 36 4 "7.271"
Set-TextCell 36 5 "  -2.47%  "
Set-TextCell 37 4 "1.031"
Set-TextCell 37 5 "  -0.77%  "
Set-TextCell 38 4 "0.2852"
Set-TextCell 38 5 "  +3.95%  "
Set-TextCell 39 4 "0.8361"
Set-TextCell 39 5 "  +5.87%  "
Set-TextCell 40 4 "0.09483"
Set-TextCell 40 5 "  +3.68%  "
Set-TextCell 41 4 "10.71"
Set-TextCell 41 5 "  -0.78%  "
Set-TextCell 42 4 "14.19"
Set-TextCell 42 5 "  -0.97%  "
Set-TextCell 43 4 "1.477"
Set-TextCell 43 5 "  +1.29%  "
Set-TextCell 44 5 "  +4.77%  "
Set-TextCell 45 4 "2.713"
Set-TextCell 45 5 "  +4.98%  "
Set-TextCell 46 4 "0.7432"
Set-TextCell 46 5 "  +2.79%  "
Set-TextCell 47 5 "  +0.39%  "
Set-TextCell 48 5 "  +3.10%  "
Set-TextCell 49 5 "  -0.21%  "
Set-TextCell 50 4 "140.28"
Set-TextCell 50 5 "  +1.64%  "
Set-TextCell 51 4 "0.08386"
Set-TextCell 51 5 "  +4.99%  "
